$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 262.9375
$ws.Range("I33").Value = 207.18182
$ws.Range("J33").Value = 385.6
$ws.Range("K33").Value = 207.18182
$ws.Range("L33").Value = 385.6
$ws.Range("M33").Value = 21.81818000000001
$ws.Range("N33").Value = -843.6

$ws.Range("H64").Value = 2602.3289
$ws.Range("I64").Value = 2514.7693
$ws.Range("J64").Value = 2694.6216
$ws.Range("K64").Value = 2514.7693
$ws.Range("L64").Value = 2694.6216
$ws.Range("M64").Value = -2266.7693
$ws.Range("N64").Value = -3190.6216

$ws.Range("H67").Value = 2602.3289
$ws.Range("I67").Value = 2514.7693
$ws.Range("J67").Value = 2694.6216
$ws.Range("K67").Value = 2514.7693
$ws.Range("L67").Value = 2694.6216
$ws.Range("M67").Value = -1656.7693
$ws.Range("N67").Value = -4410.6216

$ws.Range("H74").Value = 3329.9
$ws.Range("I74").Value = 2912.375
$ws.Range("K74").Value = 2912.375
$ws.Range("M74").Value = -1976.375

$ws.Range("H77").Value = 3329.9
$ws.Range("I77").Value = 2912.375
$ws.Range("K77").Value = 14561.875
$ws.Range("M77").Value = -9881.875

$ws.Range("H106").Value = 41669596
$ws.Range("I106").Value = 83335940
$ws.Range("J106").Value = 3258.1667
$ws.Range("K106").Value = 83335940
$ws.Range("L106").Value = 3258.1667
$ws.Range("M106").Value = -83335309
$ws.Range("N106").Value = -4520.1667

$ws.Range("H125").Value = 966.3333
$ws.Range("I125").Value = 949.5
$ws.Range("J125").Value = 1000
$ws.Range("K125").Value = 8545.5
$ws.Range("L125").Value = 9000
$ws.Range("M125").Value = -6085.5
$ws.Range("N125").Value = -13920

$ws.Range("H132").Value = 8777587
$ws.Range("I132").Value = 12200491
$ws.Range("J132").Value = 6393.5625
$ws.Range("K132").Value = 36601473
$ws.Range("L132").Value = 19180.6875
$ws.Range("M132").Value = -36598943
$ws.Range("N132").Value = -24240.6875

$ws.Range("H141").Value = 2176.0557
$ws.Range("I141").Value = 1267.8372
$ws.Range("J141").Value = 5726.364
$ws.Range("K141").Value = 3803.5116
$ws.Range("L141").Value = 17179.092
$ws.Range("M141").Value = 1376.4884
$ws.Range("N141").Value = -27539.092

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 10420485
$ws.Range("I32").Value = 3566.956
$ws.Range("J32").Value = 200008400
$ws.Range("K32").Value = 3566.956
$ws.Range("L32").Value = 200008400
$ws.Range("M32").Value = -3279.956
$ws.Range("N32").Value = -200008974

$ws.Range("H61").Value = 1079.1892
$ws.Range("I61").Value = 980.2941
$ws.Range("J61").Value = 2200
$ws.Range("K61").Value = 980.2941
$ws.Range("L61").Value = 2200
$ws.Range("M61").Value = -768.2941
$ws.Range("N61").Value = -2624

$ws.Range("H74").Value = 1257.125
$ws.Range("I74").Value = 1276.862
$ws.Range("J74").Value = 1066.3334
$ws.Range("K74").Value = 1276.862
$ws.Range("L74").Value = 1066.3334
$ws.Range("M74").Value = -402.8620000000001
$ws.Range("N74").Value = -2814.3334

$ws.Range("H77").Value = 1257.125
$ws.Range("I77").Value = 1276.862
$ws.Range("J77").Value = 1066.3334
$ws.Range("K77").Value = 6384.31
$ws.Range("L77").Value = 5331.666999999999
$ws.Range("M77").Value = -2016.31
$ws.Range("N77").Value = -14067.667

$ws.Range("H132").Value = 18577384
$ws.Range("I132").Value = 21740452
$ws.Range("J132").Value = 5350004.5
$ws.Range("K132").Value = 65221356
$ws.Range("L132").Value = 16050013.5
$ws.Range("M132").Value = -65218826
$ws.Range("N132").Value = -16055073.5

$ws.Range("H136").Value = 1079.1892
$ws.Range("I136").Value = 980.2941
$ws.Range("J136").Value = 2200
$ws.Range("K136").Value = 2940.8823
$ws.Range("L136").Value = 6600
$ws.Range("M136").Value = -390.8822999999998
$ws.Range("N136").Value = -11700

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 839
$ws.Range("I94").Value = 790.6
$ws.Range("J94").Value = 919.6667
$ws.Range("K94").Value = 790.6
$ws.Range("L94").Value = 919.6667
$ws.Range("M94").Value = -339.6
$ws.Range("N94").Value = -1821.6667

$ws.Range("H105").Value = 62501656
$ws.Range("I105").Value = 1764.9286
$ws.Range("J105").Value = 500000900
$ws.Range("K105").Value = 1764.9286
$ws.Range("L105").Value = 500000900
$ws.Range("M105").Value = -17.92859999999996
$ws.Range("N105").Value = -500004394

$ws.Range("H107").Value = 62501240
$ws.Range("I107").Value = 100000936
$ws.Range("J107").Value = 1737.6666
$ws.Range("K107").Value = 100000936
$ws.Range("L107").Value = 1737.6666
$ws.Range("M107").Value = -99999016
$ws.Range("N107").Value = -5577.6666

$ws.Range("H134").Value = 2713484.8
$ws.Range("I134").Value = 819
$ws.Range("J134").Value = 11122749
$ws.Range("K134").Value = 2457
$ws.Range("L134").Value = 33368247
$ws.Range("M134").Value = 78
$ws.Range("N134").Value = -33373317

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H132").Value = 6174072
$ws.Range("I132").Value = 1001.2973
$ws.Range("J132").Value = 19609578
$ws.Range("K132").Value = 3003.8919
$ws.Range("L132").Value = 58828734
$ws.Range("M132").Value = -473.8918999999996
$ws.Range("N132").Value = -58833794

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H116").Value = 125001870
$ws.Range("I116").Value = 10
$ws.Range("J116").Value = 142859280
$ws.Range("K116").Value = 30
$ws.Range("L116").Value = 428577840
$ws.Range("M116").Value = 3412
$ws.Range("N116").Value = -428584724

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 16671650
$ws.Range("I80").Value = 8200
$ws.Range("J80").Value = 33335100
$ws.Range("K80").Value = 8200
$ws.Range("L80").Value = 33335100
$ws.Range("M80").Value = -7202
$ws.Range("N80").Value = -33337096

$ws.Range("H83").Value = 16671650
$ws.Range("I83").Value = 8200
$ws.Range("J83").Value = 33335100
$ws.Range("K83").Value = 41000
$ws.Range("L83").Value = 166675500
$ws.Range("M83").Value = -36008
$ws.Range("N83").Value = -166685484

$ws.Range("H126").Value = 2757.4285
$ws.Range("I126").Value = 1500
$ws.Range("J126").Value = 2854.1538
$ws.Range("K126").Value = 4500
$ws.Range("L126").Value = 8562.4614
$ws.Range("M126").Value = -2030
$ws.Range("N126").Value = -13502.4614

$ws.Range("H132").Value = 6213.6387
$ws.Range("I132").Value = 3791.8125
$ws.Range("J132").Value = 25588.25
$ws.Range("K132").Value = 11375.4375
$ws.Range("L132").Value = 76764.75
$ws.Range("M132").Value = -8845.4375
$ws.Range("N132").Value = -81824.75

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 17549254
$ws.Range("I132").Value = 23257284
$ws.Range("K132").Value = 69771852
$ws.Range("M132").Value = -69769322

$ws.Range("H133").Value = 33500
$ws.Range("J133").Value = 33500
$ws.Range("L133").Value = 33500
$ws.Range("N133").Value = -38560

$ws.Range("H139").Value = 54513.5
$ws.Range("I139").Value = 53001
$ws.Range("J139").Value = 54816
$ws.Range("K139").Value = 53001
$ws.Range("L139").Value = 54816
$ws.Range("M139").Value = -47861
$ws.Range("N139").Value = -65096

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 1313.25
$ws.Range("I81").Value = 1214.8572
$ws.Range("J81").Value = 2002
$ws.Range("K81").Value = 2429.7144
$ws.Range("L81").Value = 4004
$ws.Range("M81").Value = -1368.7144
$ws.Range("N81").Value = -6126

$ws.Range("H84").Value = 1313.25
$ws.Range("I84").Value = 1214.8572
$ws.Range("J84").Value = 2002
$ws.Range("K84").Value = 12148.572
$ws.Range("L84").Value = 20020
$ws.Range("M84").Value = -6844.572
$ws.Range("N84").Value = -30628

$ws.Range("H132").Value = 21970.055
$ws.Range("I132").Value = 26584.7
$ws.Range("K132").Value = 79754.10000000001
$ws.Range("M132").Value = -77224.10000000001

Write-Output "Applied Gungnir_Profits price/profit updates across sheets."